# Append a new block of tracker rows (2025-09-16) to the progress history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 47

$goals = @(
    @{ Id = "G2"; Name = "Workout" },
    @{ Id = "G3"; Name = "Eat Healthy" },
    @{ Id = "G4"; Name = "Read Book" },
    @{ Id = "G5"; Name = "Investment Plan" },
    @{ Id = "G6"; Name = "Spend 10 Hours without phone" }
)

$dateSerial = 45916
$progress = 0.914339824239913
$percentage = 0
$change = -0.01

for ($i = 0; $i -lt $goals.Count; $i++) {
    $r = $startRow + $i
    $goal = $goals[$i]

    $ws.Cells.Item($r, 1).Value = $goal.Id
    $ws.Cells.Item($r, 2).Value = $goal.Name

    $dateCell = $ws.Cells.Item($r, 3)
    $dateCell.Value = $dateSerial
    $dateCell.NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = $progress
    $ws.Cells.Item($r, 5).Value = $percentage
    $ws.Cells.Item($r, 6).Value = $change
}
